$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data for the unit user type table (header + 6 data rows)
$ws.Range("A1").Value = "id_unit_user_type"
$ws.Range("B1").Value = "name"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "single person household"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "couple without resident child"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "couple with resident child(ren)"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "single parent with resident child(ren)"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "other household"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "company for all tertiary sectors"

# Remove the now-unused rows 8 and 9 (table shrinks from A1:B9 to A1:B7)
$ws.Range("A8:B9").ClearContents()

# Resize the table to match the new data range
$wb.Worksheets.Item(1).ListObjects.Item(1).Resize($ws.Range("A1:B7"))

# Update selection to match the post-edit state
$ws.Range("A13").Select()
